# Weekly data refresh for "Hortaliza, Terminal La Palmera de La Serena - Zapallo".
# Insert two new rows of price data at the top of the data block (row 277),
# pushing all existing rows (277..367) down by two rows to (279..369).
# The newly freed rows 277/278 are filled with the latest week's readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 277 (shifts 277..367 -> 279..369,
# inheriting each row's existing formatting, e.g. the date style on column D).
$ws.Rows.Item(277).Insert()
$ws.Rows.Item(277).Insert()

# --- New row 277: Camote "1a nueva(o)" ---
$ws.Cells.Item(277, 1).Value = 8
$ws.Cells.Item(277, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(277, 3).Value = "Coquimbo"
$ws.Cells.Item(277, 4).Value = 44468
$ws.Cells.Item(277, 5).Value = 4
$ws.Cells.Item(277, 6).Value = 100112045
$ws.Cells.Item(277, 7).Value = "Zapallo"
$ws.Cells.Item(277, 8).Value = "Camote"
$ws.Cells.Item(277, 9).Value = "1a nueva(o)"
$ws.Cells.Item(277, 10).Value = 800
$ws.Cells.Item(277, 11).Value = 700
$ws.Cells.Item(277, 12).Value = 750
$ws.Cells.Item(277, 13).Value = 725
$ws.Cells.Item(277, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(277, 15).Value = "Perú"
$ws.Cells.Item(277, 16).Value = 725
$ws.Cells.Item(277, 17).Value = 1
$ws.Cells.Item(277, 18).Value = "Hortaliza"

# --- New row 278: Camote "2a nueva(o)" ---
$ws.Cells.Item(278, 1).Value = 8
$ws.Cells.Item(278, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(278, 3).Value = "Coquimbo"
$ws.Cells.Item(278, 4).Value = 44468
$ws.Cells.Item(278, 5).Value = 4
$ws.Cells.Item(278, 6).Value = 100112045
$ws.Cells.Item(278, 7).Value = "Zapallo"
$ws.Cells.Item(278, 8).Value = "Camote"
$ws.Cells.Item(278, 9).Value = "2a nueva(o)"
$ws.Cells.Item(278, 10).Value = 520
$ws.Cells.Item(278, 11).Value = 600
$ws.Cells.Item(278, 12).Value = 650
$ws.Cells.Item(278, 13).Value = 625
$ws.Cells.Item(278, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(278, 15).Value = "Perú"
$ws.Cells.Item(278, 16).Value = 625
$ws.Cells.Item(278, 17).Value = 1
$ws.Cells.Item(278, 18).Value = "Hortaliza"
